# Updated cryptos list (price + 1h volume change) for the GitHub Actions scrape.
# Note: some "Price" values look numeric (e.g. "1.01"); a leading apostrophe
# forces Excel to keep them as text, matching the sheet's original formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.976.88'
$ws.Range('E2').Value = '  +0.50%  '
$ws.Range('D3').Value = '1.640.89'
$ws.Range('E3').Value = '  +0.92%  '
$ws.Range('E4').Value = '  +0.53%  '
$ws.Range('D5').Value = '''216.08'
$ws.Range('E5').Value = '  +0.77%  '
$ws.Range('E6').Value = '  +1.52%  '
$ws.Range('E7').Value = '  +0.46%  '
$ws.Range('E8').Value = '  +0.82%  '
$ws.Range('E9').Value = '  +1.33%  '
$ws.Range('D10').Value = '''19.64'
$ws.Range('E10').Value = '  -0.06%  '
$ws.Range('E11').Value = '  +1.12%  '
$ws.Range('E12').Value = '  +1.14%  '
$ws.Range('D13').Value = '1.866.86'
$ws.Range('E13').Value = '  +0.87%  '
$ws.Range('D14').Value = '1.658.63'
$ws.Range('E14').Value = '  +2.47%  '
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('E16').Value = '  +1.55%  '
$ws.Range('D17').Value = '''62.97'
$ws.Range('D18').Value = '25.947.99'
$ws.Range('E18').Value = '  +0.46%  '
$ws.Range('E19').Value = '  +0.51%  '
$ws.Range('D20').Value = '''193.16'
$ws.Range('E20').Value = '  +0.31%  '
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('E23').Value = '  +0.40%  '
$ws.Range('E24').Value = '  +7.50%  '
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('D26').Value = '''144.38'
$ws.Range('E26').Value = '  +2.41%  '
$ws.Range('D27').Value = '''1.01'
$ws.Range('E27').Value = '  +0.70%  '
$ws.Range('D28').Value = '''6.91'
$ws.Range('E28').Value = '  +0.99%  '
$ws.Range('D29').Value = '''15.54'
$ws.Range('E29').Value = '  +0.69%  '
$ws.Range('E30').Value = '  +0.73%  '
$ws.Range('E31').Value = '  +0.41%  '
$ws.Range('E32').Value = '  -0.93%  '
$ws.Range('E33').Value = '  +1.31%  '
$ws.Range('E34').Value = '  -2.59%  '
$ws.Range('D35').Value = '''2.46'
$ws.Range('E35').Value = '  +3.02%  '
$ws.Range('E36').Value = '  +0.34%  '
$ws.Range('D37').Value = '1.132.61'
$ws.Range('E37').Value = '  +0.43%  '
$ws.Range('E38').Value = '  -0.99%  '
$ws.Range('E39').Value = '  -0.29%  '
$ws.Range('E40').Value = '  +0.57%  '
$ws.Range('E41').Value = '  +1.56%  '
$ws.Range('D42').Value = '''99.24'
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('D43').Value = '''0.796'
$ws.Range('E43').Value = '  +0.21%  '
$ws.Range('D44').Value = '1.776.89'
$ws.Range('E44').Value = '  +0.89%  '
$ws.Range('D46').Value = '''56.57'
$ws.Range('E46').Value = '  +0.67%  '
$ws.Range('E47').Value = '  +0.79%  '
$ws.Range('D48').Value = '''1.46'
$ws.Range('E48').Value = '  +0.31%  '
$ws.Range('D49').Value = '''7.72'
$ws.Range('E49').Value = '  +1.96%  '
$ws.Range('E51').Value = '  +0.27%  '
